$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7011364102363586
$ws.Range("B1").Value = 3.484041929244995
$ws.Range("C1").Value = 2.741598129272461
$ws.Range("D1").Value = 2.268117189407349
$ws.Range("E1").Value = 2.044236660003662
